$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-11) holds the "Förändrad" (last changed) date as an
# Excel serial date number. Bump each from 45243 (2023-11-13) to
# 45244 (2023-11-14), as part of the automatic daily update.
$ws.Range("C2:C11").Value = 45244
